$d = $word.ActiveDocument

$d.Content.Find.Execute('2024-05-17 Friday', $true, $false, $false, $false, $false, $true, 1, $false, '2024-05-18 Saturday', 2) | Out-Null
$d.Content.Find.Execute('698×5=3490', $true, $false, $false, $false, $false, $true, 1, $false, '925×2=1850', 2) | Out-Null
$d.Content.Find.Execute('521×2=1042', $true, $false, $false, $false, $false, $true, 1, $false, '890×2=1780', 2) | Out-Null
$d.Content.Find.Execute('337×8=2696', $true, $false, $false, $false, $false, $true, 1, $false, '371×2=742', 2) | Out-Null
$d.Content.Find.Execute('779×7=5453', $true, $false, $false, $false, $false, $true, 1, $false, '978×7=6846', 2) | Out-Null
$d.Content.Find.Execute('319×6=1914', $true, $false, $false, $false, $false, $true, 1, $false, '471×6=2826', 2) | Out-Null
$d.Content.Find.Execute('855×9=7695', $true, $false, $false, $false, $false, $true, 1, $false, '461×3=1383', 2) | Out-Null
$d.Content.Find.Execute('506×6=3036', $true, $false, $false, $false, $false, $true, 1, $false, '112×8=896', 2) | Out-Null
$d.Content.Find.Execute('134×7=938', $true, $false, $false, $false, $false, $true, 1, $false, '837×2=1674', 2) | Out-Null
$d.Content.Find.Execute('120×5=600', $true, $false, $false, $false, $false, $true, 1, $false, '410×8=3280', 2) | Out-Null
$d.Content.Find.Execute('700×6=4200', $true, $false, $false, $false, $false, $true, 1, $false, '453×8=3624', 2) | Out-Null
$d.Content.Find.Execute('743×6=4458', $true, $false, $false, $false, $false, $true, 1, $false, '154×2=308', 2) | Out-Null
$d.Content.Find.Execute('734×8=5872', $true, $false, $false, $false, $false, $true, 1, $false, '480×8=3840', 2) | Out-Null
$d.Content.Find.Execute('344×3=1032', $true, $false, $false, $false, $false, $true, 1, $false, '660×8=5280', 2) | Out-Null
$d.Content.Find.Execute('634×7=4438', $true, $false, $false, $false, $false, $true, 1, $false, '566×4=2264', 2) | Out-Null
$d.Content.Find.Execute('663×4=2652', $true, $false, $false, $false, $false, $true, 1, $false, '639×6=3834', 2) | Out-Null
$d.Content.Find.Execute('148×9=1332', $true, $false, $false, $false, $false, $true, 1, $false, '330×2=660', 2) | Out-Null
$d.Content.Find.Execute('560×8=4480', $true, $false, $false, $false, $false, $true, 1, $false, '233×9=2097', 2) | Out-Null
$d.Content.Find.Execute('134×8=1072', $true, $false, $false, $false, $false, $true, 1, $false, '178×6=1068', 2) | Out-Null
$d.Content.Find.Execute('327×4=1308', $true, $false, $false, $false, $false, $true, 1, $false, '701×9=6309', 2) | Out-Null
$d.Content.Find.Execute('542×5=2710', $true, $false, $false, $false, $false, $true, 1, $false, '190×6=1140', 2) | Out-Null
$d.Content.Find.Execute('480×4=1920', $true, $false, $false, $false, $false, $true, 1, $false, '262×6=1572', 2) | Out-Null
$d.Content.Find.Execute('292×9=2628', $true, $false, $false, $false, $false, $true, 1, $false, '767×6=4602', 2) | Out-Null
$d.Content.Find.Execute('295×2=590', $true, $false, $false, $false, $false, $true, 1, $false, '819×5=4095', 2) | Out-Null
$d.Content.Find.Execute('807×9=7263', $true, $false, $false, $false, $false, $true, 1, $false, '910×9=8190', 2) | Out-Null
$d.Content.Find.Execute('876×4=3504', $true, $false, $false, $false, $false, $true, 1, $false, '347×2=694', 2) | Out-Null
